# aggiornamento librerie e formule
# Inserts a new "INDICATOR_36" entry into the Library_Formula sheet at row 41,
# pushing the existing rows 41-162 down to 42-163.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")
$ws.Activate()

# Insert a new row at position 41 (Excel copies formatting down from the row
# above, matching the style shift seen in the target workbook).
$ws.Rows.Item(41).Insert()

# Populate the new row with the new indicator entry.
$ws.Range("A41").Value = "CREATE/MODIFY"
$ws.Range("B41").Value = "LIB_EWS_IT"
$ws.Range("C41").Value = "INDICATOR_36"
$ws.Range("E41").Value = "String"
$ws.Range("F41").Value = "String"

# Match the recorded selection/scroll state left behind by the edit.
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G41").Select() | Out-Null
